# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# --- Sheet "Rushing" ---
$ws1 = $wb.Worksheets.Item("Rushing")

# Update existing rows 2-6
$rushingExisting = @(
    @(2, 0, "M.Stafford",   4,  3,  5,  3),
    @(3, 1, "J.Wolford",    0,  0,  1,  0),
    @(4, 2, "D.Henderson", 93, 45, 20, 27),
    @(5, 3, "S.Michel",    50, 25, 11, 20),
    @(6, 4, "J.Funk",       0,  1,  0,  0)
)
foreach ($row in $rushingExisting) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 4).Value = $row[4]
    $ws1.Cells.Item($r, 5).Value = $row[5]
    $ws1.Cells.Item($r, 6).Value = $row[6]
}

# Add new rows 7-9, copying the format of column A from row 6 (s="1")
$rushingNew = @(
    @(7, 5, "M.Brown",    0, 0, 1, 0),
    @(8, 6, "M.Sargent",  0, 1, 1, 0),
    @(9, 7, "C.Kupp",     0, 0, 1, 1)
)
foreach ($row in $rushingNew) {
    $r = $row[0]
    $ws1.Range("A6").Copy() | Out-Null
    $ws1.Range("A$r").PasteSpecial(-4122) | Out-Null
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 4).Value = $row[4]
    $ws1.Cells.Item($r, 5).Value = $row[5]
    $ws1.Cells.Item($r, 6).Value = $row[6]
}

# --- Sheet "Receiving" ---
$ws2 = $wb.Worksheets.Item("Receiving")

# Update existing rows 2-7
$receivingExisting = @(
    @(2, 0, "D.Henderson", 36, 25,  2,  2,  9,  3),
    @(3, 1, "S.Michel",    17, 14,  0,  0,  2,  1),
    @(4, 2, "C.Kupp",     103, 81, 36, 19, 25, 17),
    @(5, 3, "V.Jefferson", 50, 30, 22,  9, 13,  4),
    @(6, 4, "B.Skowronek", 17,  3,  3,  2,  1,  0),
    @(7, 5, "O.Beckham",   49, 44, 24, 10, 19, 10)
)
foreach ($row in $receivingExisting) {
    $r = $row[0]
    $ws2.Cells.Item($r, 1).Value = $row[1]
    $ws2.Cells.Item($r, 2).Value = $row[2]
    $ws2.Cells.Item($r, 3).Value = $row[3]
    $ws2.Cells.Item($r, 4).Value = $row[4]
    $ws2.Cells.Item($r, 5).Value = $row[5]
    $ws2.Cells.Item($r, 6).Value = $row[6]
    $ws2.Cells.Item($r, 7).Value = $row[7]
    $ws2.Cells.Item($r, 8).Value = $row[8]
}

# Add new rows 8-10, copying the format of column A from row 7 (s="1")
$receivingNew = @(
    @(8,  6, "K.Blanton",  2,  1, 0, 0,  1,  0),
    @(9,  7, "T.Higbee",  54, 41, 8, 3, 16, 11),
    @(10, 8, "J.Mundt",    1,  1, 0, 0,  0,  0)
)
foreach ($row in $receivingNew) {
    $r = $row[0]
    $ws2.Range("A7").Copy() | Out-Null
    $ws2.Range("A$r").PasteSpecial(-4122) | Out-Null
    $ws2.Cells.Item($r, 1).Value = $row[1]
    $ws2.Cells.Item($r, 2).Value = $row[2]
    $ws2.Cells.Item($r, 3).Value = $row[3]
    $ws2.Cells.Item($r, 4).Value = $row[4]
    $ws2.Cells.Item($r, 5).Value = $row[5]
    $ws2.Cells.Item($r, 6).Value = $row[6]
    $ws2.Cells.Item($r, 7).Value = $row[7]
    $ws2.Cells.Item($r, 8).Value = $row[8]
}
